# ProvarCache appends a newly-queried Engineering Item record to the
# "Routing Master" lookup sheet. The Item Number (col B) and Salesforce
# record Id (col D) are updated to the newly cached values; the Item
# Description (col C) already matches the cached record and is left as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Routing Master")

$ws.Range("B2").Value = "Pro-PEItem-4EFMG"
$ws.Range("D2").Value = "a345f000000uL5gAAE"

# Columns are "best fit" in the source template; re-fit them now that the
# cached values changed so the sheet keeps displaying the full text.
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(4).AutoFit()
